# Update Name of Algo
# Apply updated KNN-imputed values to columns A and C for the affected rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "A11" = -22.149
    "C11" = -13.428
    "A12" = -21.544
    "A15" = -22.062
    "C23" = -12.309
    "A27" = -21.734
    "A28" = -21.855
    "C28" = -13.14
    "A31" = -21.719
    "A32" = -22.056
    "C32" = -13.538
    "C34" = -12.345
    "A36" = -20.176
    "C36" = -12.382
    "C37" = -13.051
    "A38" = -19.722
    "C42" = -12.511
    "A46" = -21.819
    "C49" = -13.188
    "A54" = -21.937
    "C54" = -12.803
    "A55" = -22.173
    "A56" = -21.974
    "A67" = -21.577
    "A69" = -21.503
    "A72" = -21.689
    "A73" = -19.994
    "C78" = -12.52
    "C80" = -12.413
    "A83" = -21.833
    "A86" = -22.241
    "A91" = -20.967
    "A93" = -21.452
    "C97" = -11.516
    "A99" = -21.642
    "C99" = -12.418
    "C100" = -12.605
    "C101" = -12.188
    "A104" = -21.175
    "A105" = -20.454
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
